# Leakage test / Jankowski comparison samples
# Rename the leak-test SampleName values (column E) for the 2014.06.05
# rows (E2:E16) to append "Leak" -- Mix1/Mix2/3N2O/10N2O/3KCO2 become
# Mix1Leak/Mix2Leak/3N2OLeak/10N2OLeak/3KCO2Leak. Rows 17:25 (amb/NA/
# 2014.06.06 samples) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Mix1Leak"
$ws.Range("E3").Value = "Mix1Leak"
$ws.Range("E4").Value = "Mix1Leak"

$ws.Range("E5").Value = "Mix2Leak"
$ws.Range("E6").Value = "Mix2Leak"
$ws.Range("E7").Value = "Mix2Leak"

$ws.Range("E8").Value = "3N2OLeak"
$ws.Range("E9").Value = "3N2OLeak"
$ws.Range("E10").Value = "3N2OLeak"

$ws.Range("E11").Value = "10N2OLeak"
$ws.Range("E12").Value = "10N2OLeak"
$ws.Range("E13").Value = "10N2OLeak"

$ws.Range("E14").Value = "3KCO2Leak"
$ws.Range("E15").Value = "3KCO2Leak"
$ws.Range("E16").Value = "3KCO2Leak"

# Leave the active selection on E3, matching where work was in progress.
$ws.Range("E3").Select()
